$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tp/fp/tn/fn values for rows 6-12 were swapped (C,D,E,F) -> (F,E,D,C),
# and the derived precision/recall/fscore columns (G,H,I) are recomputed
# from the corrected tp/fp/fn values.
$rows = @(
    @{ Row = 6;  C = 135; D = 16; E = 0; F = 0 },
    @{ Row = 7;  C = 135; D = 16; E = 0; F = 0 },
    @{ Row = 8;  C = 135; D = 16; E = 0; F = 0 },
    @{ Row = 9;  C = 135; D = 16; E = 0; F = 0 },
    @{ Row = 10; C = 132; D = 16; E = 0; F = 0 },
    @{ Row = 11; C = 130; D = 16; E = 0; F = 0 },
    @{ Row = 12; C = 121; D = 14; E = 0; F = 0 }
)

foreach ($r in $rows) {
    $tp = $r.C
    $fp = $r.D
    $tn = $r.E
    $fn = $r.F

    if (($tp + $fp) -ne 0) {
        $precision = $tp / ($tp + $fp)
    } else {
        $precision = 0
    }

    if (($tp + $fn) -ne 0) {
        $recall = $tp / ($tp + $fn)
    } else {
        $recall = 0
    }

    if (($precision + $recall) -ne 0) {
        $fscore = 2 * $precision * $recall / ($precision + $recall)
    } else {
        $fscore = 0
    }

    $ws.Cells.Item($r.Row, 3).Value = $tp
    $ws.Cells.Item($r.Row, 4).Value = $fp
    $ws.Cells.Item($r.Row, 5).Value = $tn
    $ws.Cells.Item($r.Row, 6).Value = $fn
    $ws.Cells.Item($r.Row, 7).Value = $precision
    $ws.Cells.Item($r.Row, 8).Value = $recall
    $ws.Cells.Item($r.Row, 9).Value = $fscore
}
